# Bangalore CMP 2020 projects workbook - add a "budget" column.
#
# The sheet currently has columns A..F:
#   A id, B type, C description_0, D description_1, E description_2, F location
# We insert a new column F ("budget") ahead of the existing location column,
# which shifts the old F (location) to G. Excel shifts all the existing data
# (incl. shared-string cell references) automatically when a whole column is
# inserted, so only the brand-new "budget" values and the handful of brand
# new location strings need to be written explicitly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert the new column; old F (location) becomes G automatically, along
#    with every existing cell value/shared-string reference in that column.
$ws.Columns("F:F").Insert()

# 2) Header for the new column.
$ws.Range("F1").Value = "budget"

# 3) Budget figures for the two project blocks that received them.
#    (ids 1001-1051 -> rows 2-52, all budget 10)
$ws.Range("F2:F52").Value = 10
#    (ids 3001-3058 -> rows 208-265, all budget 50)
$ws.Range("F208:F265").Value = 50

# 4) A few rows at the bottom of the sheet gained brand-new location values
#    in column G (they had no prior "location" data in column F at all).
#    Write them in ascending shared-string-allocation order to match the
#    authoring tool's order (budget already consumed the first new slot).
$ws.Range("G311").Value = "13.013255992883385, 77.76103920899732"
$ws.Range("G310").Value = "12.914265937865533, 77.48636936233684"
$ws.Range("G309").Value = "12.916671511371238, 77.48251674160508"
$ws.Range("G305").Value = "13.018173812026925, 77.55657081141072"

# 5) Column widths: new F (numeric budget column) matches the narrow
#    "13.71" width used by the other numeric columns (D/E); G keeps the
#    location column's original width, which is preserved automatically by
#    the column insert above.
$ws.Columns("F:F").ColumnWidth = $ws.Columns("D:D").ColumnWidth()

# 6) Restore the frozen-pane selection to the bottom pane's new active cell.
$ws.Range("B1").Select()
$ws.Range("E269").Select()
